# Add a new "Ethnic or Racial Identity" field to the SA-MODS template.
# This mirrors the existing "Gender Identity, Sexuality" triplet of columns
# (G:I -> <mods:description>…: / label / </mods:description>) by inserting
# three new columns right after it (at J:L) and filling them in the same
# pattern, which pushes every subsequent column three places to the right.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three blank columns before the old column J (Pronouns' opening tag),
# shifting everything from J onward to the right by 3.
$ws.Range("J1:L1").EntireColumn.Insert() | Out-Null

# Populate the new triplet of cells with the new field's markup/label,
# following the same structure used by the other description fields.
$ws.Range("J1").Value = "<mods:description>Ethnic or Racial Identity: "
$ws.Range("K1").Value = "Ethnic or Racial Identity"
$ws.Range("L1").Value = "</mods:description>"

# The label cell (middle of the triplet) uses the bold header style, same
# as the other field-name cells (e.g. H1 "Gender Identity, Sexuality").
$ws.Range("K1").Font.Bold = $true

# Excel left the active cell on the newly inserted column after editing.
$ws.Range("J1").Select() | Out-Null
